$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row 7 (shifts nothing currently there; extends the table by one weekly record)
$ws.Rows.Item(7).Insert()

# Row 7 receives the data that used to be in row 6 before this edit
$ws.Cells.Item(7,1).Value = 7
$ws.Cells.Item(7,2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(7,3).Value = "Ñuble"
$ws.Cells.Item(7,4).Value = 44608
$ws.Cells.Item(7,4).NumberFormat = $ws.Cells.Item(6,4).NumberFormat
$ws.Cells.Item(7,5).Value = 16
$ws.Cells.Item(7,6).Value = 100112044
$ws.Cells.Item(7,7).Value = "Perejil"
$ws.Cells.Item(7,8).Value = "Sin especificar"
$ws.Cells.Item(7,9).Value = "Primera"
$ws.Cells.Item(7,10).Value = 120
$ws.Cells.Item(7,11).Value = 600
$ws.Cells.Item(7,12).Value = 650
$ws.Cells.Item(7,13).Value = 625
$ws.Cells.Item(7,14).Value = "$/atado 0,5 a 1 kilo"
$ws.Cells.Item(7,15).Value = "Región del Maule"
$ws.Cells.Item(7,16).Value = 625
$ws.Cells.Item(7,17).Value = 1
$ws.Cells.Item(7,18).Value = "Hortaliza"

# Row 6 is updated with the values that used to be in row 5 before this edit
$ws.Cells.Item(6,4).Value = 44624
$ws.Cells.Item(6,11).Value = 650
$ws.Cells.Item(6,12).Value = 700
$ws.Cells.Item(6,13).Value = 675
$ws.Cells.Item(6,16).Value = 675

# Row 5 is updated with the new weekly record
$ws.Cells.Item(5,4).Value = 44754
$ws.Cells.Item(5,10).Value = 200
$ws.Cells.Item(5,11).Value = 700
$ws.Cells.Item(5,12).Value = 750
$ws.Cells.Item(5,13).Value = 725
$ws.Cells.Item(5,16).Value = 725
